$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor expression stats and derived edge weights/specificities
# (new TPM-based recomputation), applied per the commit diff.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.176022
$ws.Range("H2").Value = 0.5280659999999999
$ws.Range("I2").Value = 0.03293066697281707
$ws.Range("J2").Value = 0.03293066697281707
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 0.4516425869339999
$ws.Range("R2").Value = 4.064783282405999
$ws.Range("S2").Value = 0.003076335519222081
$ws.Range("T2").Value = 0.003076335519222081
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.176022
$ws.Range("H3").Value = 0.5280659999999999
$ws.Range("I3").Value = 0.03293066697281707
$ws.Range("J3").Value = 0.03293066697281707
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 1.860261458286
$ws.Range("R3").Value = 16.742353124574
$ws.Range("S3").Value = 0.01267105575232519
$ws.Range("T3").Value = 0.01267105575232518
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.176022
$ws.Range("H4").Value = 0.5280659999999999
$ws.Range("I4").Value = 0.03293066697281707
$ws.Range("J4").Value = 0.03293066697281707
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 2.408643037416
$ws.Range("R4").Value = 21.677787336744
$ws.Range("S4").Value = 0.01640632292767515
$ws.Range("T4").Value = 0.01640632292767515
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.176022
$ws.Range("H5").Value = 0.5280659999999999
$ws.Range("I5").Value = 0.03293066697281707
$ws.Range("J5").Value = 0.03293066697281707
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 0.114065893788
$ws.Range("R5").Value = 1.026593044092
$ws.Range("S5").Value = 0.0007769527735946581
$ws.Range("T5").Value = 0.000776952773594658
# Row 6
$ws.Range("I6").Value = 0.8002039325901205
$ws.Range("J6").Value = 0.8002039325901203
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 10.97476022845478
$ws.Range("R6").Value = 98.772842056093
$ws.Range("S6").Value = 0.07475389983689699
$ws.Range("T6").Value = 0.07475389983689698
# Row 7
$ws.Range("I7").Value = 0.8002039325901205
$ws.Range("J7").Value = 0.8002039325901203
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("S7").Value = 0.3079023164471271
$ws.Range("T7").Value = 0.3079023164471271
# Row 8
$ws.Range("I8").Value = 0.8002039325901205
$ws.Range("J8").Value = 0.8002039325901203
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 58.52920113452578
$ws.Range("R8").Value = 526.762810210732
$ws.Range("S8").Value = 0.3986680299219593
$ws.Range("T8").Value = 0.3986680299219593
# Row 9
$ws.Range("I9").Value = 0.8002039325901205
$ws.Range("J9").Value = 0.8002039325901203
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 2.771762165002889
$ws.Range("R9").Value = 24.945859485026
$ws.Range("S9").Value = 0.01887968638413707
$ws.Range("T9").Value = 0.01887968638413707
# Row 10
$ws.Range("G10").Value = 0.891934
$ws.Range("H10").Value = 2.675802
$ws.Range("I10").Value = 0.1668654004370625
$ws.Range("J10").Value = 0.1668654004370625
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 2.288551312531333
$ws.Range("R10").Value = 20.596961812782
$ws.Range("S10").Value = 0.01558832557863124
$ws.Range("T10").Value = 0.01558832557863124
# Row 11
$ws.Range("G11").Value = 0.891934
$ws.Range("H11").Value = 2.675802
$ws.Range("I11").Value = 0.1668654004370625
$ws.Range("J11").Value = 0.1668654004370625
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 9.426267418475332
$ws.Range("R11").Value = 84.836406766278
$ws.Range("S11").Value = 0.0642064369305792
$ws.Range("T11").Value = 0.0642064369305792
# Row 12
$ws.Range("G12").Value = 0.891934
$ws.Range("H12").Value = 2.675802
$ws.Range("I12").Value = 0.1668654004370625
$ws.Range("J12").Value = 0.1668654004370625
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 12.20501198108533
$ws.Range("R12").Value = 109.845107829768
$ws.Range("S12").Value = 0.08313368348372936
$ws.Range("T12").Value = 0.08313368348372936
# Row 13
$ws.Range("G13").Value = 0.891934
$ws.Range("H13").Value = 2.675802
$ws.Range("I13").Value = 0.1668654004370625
$ws.Range("J13").Value = 0.1668654004370625
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 0.5779916653026667
$ws.Range("R13").Value = 5.201924987723999
$ws.Range("S13").Value = 0.003936954444122768
$ws.Range("T13").Value = 0.003936954444122768
